$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text formatting so values
# like "1.00", "99.20", "37.815.98" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.815.98"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "2.086.06"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "233.56"
$ws.Range("E5").Value = "  +0.46%  "

$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("D7").Value = "58.87"
$ws.Range("E7").Value = "  +3.58%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("D10").Value = "0.0791"
$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("E11").Value = "  +2.95%  "

$ws.Range("D12").Value = "2.391.76"
$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").Value = "14.74"
$ws.Range("E13").Value = "  +2.63%  "

$ws.Range("D14").Value = "21.27"
$ws.Range("E14").Value = "  +2.19%  "

$ws.Range("D15").Value = "0.777"
$ws.Range("E15").Value = "  +2.83%  "

$ws.Range("D16").Value = "5.34"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("D17").Value = "2.078.78"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").Value = "37.745.43"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").Value = "6.18"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").Value = "71.87"
$ws.Range("E20").Value = "  +1.83%  "

$ws.Range("E21").Value = "  +3.60%  "

$ws.Range("D22").Value = "228.44"
$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("D26").Value = "9.57"
$ws.Range("E26").Value = "  +7.86%  "

$ws.Range("D27").Value = "171.31"
$ws.Range("E27").Value = "  +0.95%  "

$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").Value = "1.44"
$ws.Range("E29").Value = "  +1.60%  "

$ws.Range("E30").Value = "  +1.38%  "

$ws.Range("E31").Value = "  +2.88%  "

$ws.Range("D32").Value = "4.75"
$ws.Range("E32").Value = "  +3.27%  "

$ws.Range("D33").Value = "0.0636"
$ws.Range("E33").Value = "  +2.42%  "

$ws.Range("D34").Value = "4.69"
$ws.Range("E34").Value = "  +2.41%  "

$ws.Range("E35").Value = "  +1.00%  "

$ws.Range("D36").Value = "3.44"
$ws.Range("E36").Value = "  +2.77%  "

$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "5.44"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("D40").Value = "0.0982"
$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "99.20"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0220"
$ws.Range("E42").Value = "  +3.30%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "17.17"
$ws.Range("E43").Value = "  +10.73%  "

$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "1.452.95"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("E47").Value = "  -4.58%  "

$ws.Range("E48").Value = "  +2.01%  "

$ws.Range("D49").Value = "7.38"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").Value = "2.277.75"
$ws.Range("E51").Value = "  +0.70%  "
